# Apply crypto price/volume update from GitHub Actions scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding numeric-looking text (e.g. "304.81") must keep their
# Text format so Excel does not coerce them into real numbers --
# the source sheet stores every Price/Volume value as a string.
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.438.64"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.631.30"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "304.81"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("D7").Value = "0.3725"
$ws.Range("E7").Value = "  -0.88%  "
$ws.Range("D8").Value = "0.3645"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "51.61"
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("D10").Value = "0.08174"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("E11").Value = "  -4.38%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "22.47"
$ws.Range("E13").Value = "  -2.43%  "
$ws.Range("D14").Value = "6.539"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").Value = "0.00001247"
$ws.Range("E15").Value = "  -2.78%  "
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").Value = "1.631.69"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "94.49"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").Value = "0.06955"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").Value = "17.75"
$ws.Range("E20").Value = "  -3.08%  "
$ws.Range("D21").Value = "6.457"
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").Value = "23.434.16"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").Value = "3.210"
$ws.Range("E25").Value = "  +4.46%  "
$ws.Range("D26").Value = "2.466"
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").Value = "21.35"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "149.68"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").Value = "5.342"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "134.38"
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("D31").Value = "1.812.34"
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("D32").Value = "2.281"
$ws.Range("E32").Value = "  -4.04%  "
$ws.Range("D33").Value = "6.821"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").Value = "1.020"
$ws.Range("E34").Value = "  +4.35%  "
$ws.Range("D35").Value = "11.04"
$ws.Range("E35").Value = "  +6.14%  "
$ws.Range("D36").Value = "0.02779"
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("D37").Value = "0.2523"
$ws.Range("E37").Value = "  -1.37%  "
$ws.Range("D38").Value = "0.08743"
$ws.Range("E38").Value = "  -1.72%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "6.051"
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.07120"
$ws.Range("E40").Value = "  -3.88%  "
$ws.Range("D41").Value = "0.7038"
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "16.33"
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "1.345"
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("D44").Value = "12.28"
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("D45").Value = "0.6511"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("D46").Value = "2.329"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").Value = "0.9999"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "3.991"
$ws.Range("E48").Value = "  -1.35%  "
$ws.Range("D49").Value = "0.08031"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").Value = "1.205"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("D51").Value = "125.03"
$ws.Range("E51").Value = "  -3.83%  "
